# Update countries & provincias Spain
# Applies:
#  - Reordering of 3 country-name entries (Rusia, Marruecos, Uzbekistan),
#    which shifts the displayed country + stats for the rows in between.
#  - Refreshed stats for a handful of rows (12, 24) unrelated to the reorder.
#  - Updated "last refreshed" timestamp string in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "datos actualizados" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 09:55"

# --- Helper: write a full data row (Pais + 7 numeric stats) ---
function Set-Row($r, $pais, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($r, 1).Value = $pais
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = $f
    $ws.Cells.Item($r, 7).Value = $g
    $ws.Cells.Item($r, 8).Value = $h
}

# --- Row 12 (Suiza): refreshed totals, no reordering involved ---
Set-Row 12 "Suiza" 17781 13 2967 14326 348 0 488

# --- Row 24 (Noruega): refreshed totals, no reordering involved ---
Set-Row 24 "Noruega" 4879 2 13 4822 105 0 44

# --- Rows 26-30: "Rusia" moves up next to "Chequia" in the country list,
#     pushing Irlanda/Dinamarca/Chile/Malasia down one row each, and
#     Rusia itself receives freshly updated stats. ---
Set-Row 26 "Rusia"     3548 771 235 3283 8   6 30
Set-Row 27 "Irlanda"   3447 0   5   3357 103 0 85
Set-Row 28 "Dinamarca" 3107 0   894 2109 145 0 104
Set-Row 29 "Chile"     3031 0   234 2781 31  0 16
Set-Row 30 "Malasia"   2908 0   645 2218 102 0 45

# --- Rows 66-67: "Marruecos" moves ahead of "Armenia", Marruecos gets
#     freshly updated stats, Armenia keeps its old stats one row down. ---
Set-Row 66 "Marruecos" 676 22 29 608 1  0 39
Set-Row 67 "Armenia"   663 92 33 626 30 0 4

# --- Rows 100-102: "Uzbekistan" moves ahead of "Senegal", Uzbekistan gets
#     freshly updated stats, Senegal/Malta shift down one row each. ---
Set-Row 100 "Uzbekistan" 190 9 12 176 8 0 2
Set-Row 101 "Senegal"    190 0 45 144 0 0 1
Set-Row 102 "Malta"      188 0 2  186 2 0 0
